# Update "想去人数" (number of people interested) values on two sheets:
# "展览" (rows 3-6) and "全部类型" (rows 4-7), column F.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 290
$wsExhibit.Range("F4").Value = 1277
$wsExhibit.Range("F5").Value = 80
$wsExhibit.Range("F6").Value = 59

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 290
$wsAll.Range("F5").Value = 1277
$wsAll.Range("F6").Value = 80
$wsAll.Range("F7").Value = 59
